# Newcastle Utd_stats.xlsx edit script
# 1) Rename the per-category stat sheets to their "spaced out" display names.
# 2) Bump every player's "Age" column (col E, format "YY-DDD") forward by one
#    day on every stats sheet (StandardStats .. MiscStats), since the stats
#    export date moved forward a day.

$wb = $excel.ActiveWorkbook

# --- 1) Sheet renames -------------------------------------------------
$renames = @{
    "StandardStats"    = "Standard Stats"
    "ShootingStats"     = "Shooting Stats"
    "PassingStats"      = "Passing Stats"
    "PassTypes"         = "Pass Types"
    "GoalShotCreation"  = "Goal & Shot Creation"
    "DefensiveActions"  = "Defensive Actions"
    "PlayingTime"       = "Playing Time"
    "MiscStats"         = "Miscellaneous Stats"
}

# Sheet names as they exist in the workbook BEFORE renaming, in order, so we
# can still find each sheet after earlier renames happen.
$origOrder = @("StandardStats","ShootingStats","PassingStats","PassTypes","GoalShotCreation","DefensiveActions","PlayingTime","MiscStats")

foreach ($oldName in $origOrder) {
    $newName = $renames[$oldName]
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $newName
}

# --- 2) Increment the "Age" (YY-DDD) values by one day ----------------
function Bump-Age([string]$ageStr) {
    if ($ageStr -match '^(\d+)-(\d+)$') {
        $year = $matches[1]
        $day = [int]$matches[2]
        $day = $day + 1
        $dayStr = ([string]$day).PadLeft(3, '0')
        return "$year-$dayStr"
    }
    return $ageStr
}

# All per-category stats sheets (their current/new names) get the Age bump.
# "Matches" has no Age column and is skipped; "Possession" wasn't renamed
# but still carries the same player Age column as the others.
$statSheetNames = @("Standard Stats","Shooting Stats","Passing Stats","Pass Types","Goal & Shot Creation","Defensive Actions","Possession","Playing Time","Miscellaneous Stats")

foreach ($sheetName in $statSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lastRow = $ws.UsedRange.Rows.Count
    for ($row = 4; $row -le $lastRow; $row++) {
        $cell = $ws.Cells.Item($row, 5)
        $val = $cell.Value2
        if ($val -ne $null) {
            $strVal = [string]$val
            if ($strVal -match '^\d+-\d+$') {
                $cell.Value2 = Bump-Age $strVal
            }
        }
    }
}
